$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42600.786805555559
$ws.Range("B3").Value = "Noun"
$ws.Range("C3").Value = 5472
$ws.Range("D3").Value = 2026
$ws.Range("E3").Value = 356
$ws.Range("F3").Value = 62
$ws.Range("G3").Value = 13
$ws.Range("H3").Value = 82
$ws.Range("I3").Value = 17
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0

$ws.Range("A4").Value = 42600.825462962966
$ws.Range("B4").Value = "Noun"
$ws.Range("C4").Value = 4727
$ws.Range("D4").Value = 1603
$ws.Range("E4").Value = 318
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 83
$ws.Range("I4").Value = 16
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0

$ws.Columns.Item(1).ColumnWidth = 14.022135416666666
